# Insert a new row at row 315 (this shifts the existing rows 315-418 down to 316-419,
# preserving their data and formatting) and populate the new row with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 315, shifting rows 315-418 -> 316-419.
$ws.Rows.Item(315).Insert()

# Populate the new row 315 with the new data record.
$ws.Cells.Item(315, 1).Value = 5
$ws.Cells.Item(315, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(315, 3).Value = "Maule"
$ws.Cells.Item(315, 4).Value = 44876
$ws.Cells.Item(315, 5).Value = 7
$ws.Cells.Item(315, 6).Value = 100114013
$ws.Cells.Item(315, 7).Value = "Zanahoria"
$ws.Cells.Item(315, 8).Value = "Sin especificar"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 500
$ws.Cells.Item(315, 11).Value = 13000
$ws.Cells.Item(315, 12).Value = 13000
$ws.Cells.Item(315, 13).Value = 13000
$ws.Cells.Item(315, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(315, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(315, 16).Value = 650
$ws.Cells.Item(315, 17).Value = 20
$ws.Cells.Item(315, 18).Value = "Hortaliza"
